# Insert a new daily price record at row 27 (the sheet is sorted/accumulated
# with the newest record inserted above the previous one for that date
# series). This pushes the former rows 27-132 down to 28-133 and extends the
# used range from A1:R132 to A1:R133.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 27, shifting rows 27:132 down to 28:133.
$ws.Rows("27:27").Insert()

# The new blank row 27 should start as a duplicate of the record that is now
# in row 28 (the original row 27), then we just update its date.
$ws.Rows("28:28").Copy()
$ws.Rows("27:27").PasteSpecial()

# Set the new record's date (column D) to the new reporting date.
$ws.Range("D27").Value = 45099
